$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:D19")
$key = $ws.Range("A2:A19")

$rng.Sort($key, 1, $null, $null, 1, $null, 1, 1, $false, $null, $null, 1)
